# Remove the four user records that were dropped from the dataset:
#   BTN0003694 - RUSTY ARI MEI MANALU   (row 4)
#   BTN0004670 - EKO RESTIOWATI         (row 5)
#   BTN0004015 - Frisilia               (row 10, before any deletion)
#   BTN0000041 - Rara                   (row 19, before any deletion)
# Deleting from the bottom up keeps the remaining row numbers stable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(19).Delete()
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# Reflect the author's final selection: the (now empty-of-those-rows)
# row 4 is selected as a whole row.
$ws.Rows.Item(4).EntireRow.Select() | Out-Null
